$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add missing X6 / Y6 values to existing row 6 ---
$ws.Range("X6").Value = 0.29000100000000373
$ws.Range("Y6").Value = "Up"

# --- Add new row 7, copying number formats from row 6 where needed ---

# A7: date/time formatted like A6
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 42648.890567129631

$ws.Range("B7").Value = -9
$ws.Range("C7").Value = "Sell"
$ws.Range("D7").Value = 18
$ws.Range("E7").Value = 15905
$ws.Range("F7").Value = 2949
$ws.Range("G7").Value = 61
$ws.Range("H7").Value = 37
$ws.Range("I7").Value = 81
$ws.Range("J7").Value = 18
$ws.Range("K7").Value = 39880
$ws.Range("L7").Value = 341
$ws.Range("M7").Value = 211
$ws.Range("N7").Value = 123
$ws.Range("O7").Value = 28
$ws.Range("P7").Value = "Noun"
$ws.Range("Q7").Value = 38.916275631518758
$ws.Range("R7").Value = 0

# S7 / T7: percent formatted like S6 / T6
$ws.Range("S6").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = -0.0078

$ws.Range("T6").Copy()
$ws.Range("T7").PasteSpecial(-4122)
$ws.Range("T7").Value = -0.0305

$ws.Range("U7").Value = 14.62
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = -2
